$d = $word.ActiveDocument

# --- Change 1: the lone empty (non-list) paragraph after "Dkien win..." switches
#     its paragraph-mark language from en-US to vi-VN ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Length -eq 1 -and $p.Style.NameLocal -eq "Normal") {
        $p.Range.LanguageID = "vi-VN"
        break
    }
}

# --- Change 2: append the new "Lvl up / LVL1 / LVL2 / Lvl3" level-tuning block
#     at the end of the document (new paragraphs with their own text/runs) ---
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$fragment = @'
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Lvl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t xml:space="preserve"> up: </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>tỉ</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t xml:space="preserve"> lệ số 2,3 tăng</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>LVL</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t>1: tổng số 2, 3 dao động từ 3-4-5. Số boom 6</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>LVL</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t>2: 4 5 6. Số boom là 7</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Lvl</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t>3: 5,6,7. Số boom laf 8</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="ListParagraph"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="2"/>
        </w:numPr>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$endRng.InsertXML($fragment)
